# Daily attendance processing - 2026-01-25 21:33:29
#
# Normalises the "Recorded By" column (G) on the "Session Analysis Results"
# sheet: each cell holds a comma-separated list of recorder names/emails
# (e.g. "dnasr281@gmail.com, System"). Re-sort each list in ascending,
# case-sensitive (ordinal) order so the values are consistently ordered,
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value.ToString().Contains(",")) {

        $rawParts = $value.ToString().Split(",")
        $n = $rawParts.Length

        $parts = @()
        for ($i = 0; $i -lt $n; $i++) {
            $parts += $rawParts[$i].Trim()
        }

        # Ordinal (case-sensitive) ascending bubble sort - Sort-Object's
        # default/-CaseSensitive comparer is culture-aware, not ordinal, so
        # it can't reproduce Excel's straightforward ordinal ordering here.
        for ($i = 0; $i -lt $n; $i++) {
            for ($j = 0; $j -lt ($n - $i - 1); $j++) {
                $cmp = $parts[$j].CompareTo($parts[$j + 1])
                if ($cmp -gt 0) {
                    $tmp = $parts[$j]
                    $parts[$j] = $parts[$j + 1]
                    $parts[$j + 1] = $tmp
                }
            }
        }

        $newValue = [string]::Join(", ", $parts)

        if ($newValue -ne $value.ToString()) {
            $cell.Value = $newValue
        }
    }
}
